$wb = $excel.ActiveWorkbook

# --- Update the Date field on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-25T15:25:45+00:00"

# --- Add a new concept row ("TNEBA") on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Copy the formatting (borders/style) of the last existing data row down
# onto the new row so the new cells carry the same cell style.
$concepts.Range("A7:D7").Copy()
$concepts.Range("A8:D8").PasteSpecial(-4122)

# Column A holds the textual level ("1"), same as every other data row.
# Copying the *value* (not the format) from an existing text-typed cell
# keeps it stored as text/shared-string rather than a numeric literal.
$concepts.Range("A2").Copy()
$concepts.Range("A8").PasteSpecial(-4163)

$concepts.Range("B8").Value = "TNEBA"
$concepts.Range("C8").Value = "Tumor-Normal Exomes Bioinformatic Analysis"
